$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: I ("I0") and J ("IF"). Add the two header cells first, then
# copy the existing "IP" header's formatting (bold font, border, centered
# alignment) onto them so they match the rest of row 1 exactly.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: column I is a constant 1 for every existing data row, and
# column J mirrors column H (IP) for that row -- except the last row (35)
# which carries its own distinct values.
$jByRow = @{
    2 = 8; 3 = 2; 4 = 9; 5 = 5; 6 = 4; 7 = 7; 8 = 7; 9 = 4; 10 = 5;
    11 = 6; 12 = 7; 13 = 6; 14 = 5; 15 = 3; 16 = 7; 17 = 7; 18 = 7;
    19 = 7; 20 = 7; 21 = 6; 22 = 5; 23 = 6; 24 = 5; 25 = 7; 26 = 7;
    27 = 6; 28 = 8; 29 = 6; 30 = 5; 31 = 5; 32 = 6; 33 = 5; 34 = 5
}

foreach ($row in 2..34) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $jByRow[$row]
}

# Row 35 is the exception: I35 = 4, J35 = 6.
$ws.Cells.Item(35, 9).Value = 4
$ws.Cells.Item(35, 10).Value = 6
